$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# The "no carga resoluciones" bug is fixed: two pending resolution
# requests ("pato" / "patito") that previously failed to load now
# show up as new rows in the submissions log.
# -----------------------------------------------------------------

# ---- Row 5: pato ----
$ws.Cells.Item(5,1).Value = 42903.51540038195
$ws.Cells.Item(5,1).NumberFormat = "d/m/yy"

$ws.Cells.Item(5,2).NumberFormat = "@"
$ws.Cells.Item(5,2).Value = "2"
$ws.Cells.Item(5,2).Style = "Normal"

$ws.Cells.Item(5,3).Value = "pato"
$ws.Cells.Item(5,4).Value = "pato@gmail.com"

$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = "60990312"
$ws.Cells.Item(5,5).Style = "Normal"

$ws.Cells.Item(5,6).Value = "IS2017"
$ws.Cells.Item(5,7).Value = "IC1802"
$ws.Cells.Item(5,8).Value = 1
$ws.Cells.Item(5,9).Value = "MODIFICACIÓN_ACTA"
$ws.Cells.Item(5,10).Value = "prueba"

$ws.Cells.Item(5,11).NumberFormat = "@"
$ws.Cells.Item(5,11).Value = "2"
$ws.Cells.Item(5,11).Style = "Normal"

$ws.Cells.Item(5,12).Value = "pato"
$ws.Cells.Item(5,13).Value = "PENDIENTE"

# ---- Row 6: patito ----
$ws.Cells.Item(6,1).Value = 42903.51825075231
$ws.Cells.Item(6,1).NumberFormat = "d/m/yy"

$ws.Cells.Item(6,2).NumberFormat = "@"
$ws.Cells.Item(6,2).Value = "2"
$ws.Cells.Item(6,2).Style = "Normal"

$ws.Cells.Item(6,3).Value = "patito"
$ws.Cells.Item(6,4).Value = "b@g.com"

$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = "20998766"
$ws.Cells.Item(6,5).Style = "Normal"

$ws.Cells.Item(6,6).Value = "IS2017"
$ws.Cells.Item(6,7).Value = "IC1802"
$ws.Cells.Item(6,8).Value = 1
$ws.Cells.Item(6,9).Value = "MODIFICACIÓN_ACTA"
$ws.Cells.Item(6,10).Value = ":V "

$ws.Cells.Item(6,11).NumberFormat = "@"
$ws.Cells.Item(6,11).Value = "2"
$ws.Cells.Item(6,11).Style = "Normal"

$ws.Cells.Item(6,12).Value = "patito"
$ws.Cells.Item(6,13).Value = "PENDIENTE"
